$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'306.76"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.35%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'38.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'7.04%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.111"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'1.28%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08089"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-0.66%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.937"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-6.78%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'7.977"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'1.37%"
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'0.19%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1482"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'1.59%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1929"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'0.11%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.09173"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'0.62%"
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'1.53%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.09785"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.87%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001406"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.08%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.005969"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-3.89%"
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'-1.51%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'4.187"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.75%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.426"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'0.76%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3459"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-0.12%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1302"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-1.25%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'4.684"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-2.95%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.2416"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'3.38%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04378"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.08%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'0.52%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004280"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'1.95%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'0.18%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D39").Value = "'0.02039"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-0.41%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.05068"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-1.83%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007541"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'1.18%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.01017"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'0.62%"
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'-1.78%"
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'-0.29%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.009921"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'2.04%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006183"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-1.74%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'0.37%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.003111"
$ws.Range("D48").Style = "Normal"
$ws.Range("E49").Value = "'0.31%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002103"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.37%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002003"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.37%"
$ws.Range("E51").Style = "Normal"
